$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = @'
2508.05619v1
'@

$ws.Range("B2").Value = @'
The Missing Reward: Active Inference in the Era of Experience
'@

$ws.Range("C2").Value = @'
Bo Wen
'@

$ws.Range("D2").Value = @'
http://arxiv.org/abs/2508.05619v1
'@

$ws.Range("E2").Value = @'
  This paper argues that Active Inference (AIF) provides a crucial foundation
for developing autonomous AI agents capable of learning from experience without
continuous human reward engineering. As AI systems begin to exhaust
high-quality training data and rely on increasingly large human workforces for
reward design, the current paradigm faces significant scalability challenges
that could impede progress toward genuinely autonomous intelligence. The
proposal for an ``Era of Experience,'' where agents learn from self-generated
data, is a promising step forward. However, this vision still depends on
extensive human engineering of reward functions, effectively shifting the
bottleneck from data curation to reward curation. This highlights what we
identify as the \textbf{grounded-agency gap}: the inability of contemporary AI
systems to autonomously formulate, adapt, and pursue objectives in response to
changing circumstances. We propose that AIF can bridge this gap by replacing
external reward signals with an intrinsic drive to minimize free energy,
allowing agents to naturally balance exploration and exploitation through a
unified Bayesian objective. By integrating Large Language Models as generative
world models with AIF's principled decision-making framework, we can create
agents that learn efficiently from experience while remaining aligned with
human values. This synthesis offers a compelling path toward AI systems that
can develop autonomously while adhering to both computational and physical
constraints.

'@

$ws.Range("F2").Value = ""

$ws.Range("I2").Value = @'
- Active Inference (AIF),   - Large Language Models (as generative world models),   - Bayesian objective for balancing exploration and exploitation,   - Free energy minimization, , -
'@

$ws.Range("J2").Value = @'
- AIF provides a foundation for developing autonomous AI agents capable of learning from experience without continuous human reward engineering.,   - AIF can bridge the grounded-agency gap by replacing external reward signals with an intrinsic drive to minimize free energy.,   - Integration of Large Language Models with AIF's decision-making framework can create agents that learn efficiently from experience while remaining aligned with human values.,   - This approach offers a path toward AI systems that can develop autonomously while adhering to computational and physical constraints., , -
'@

$ws.Range("K2").Value = @'
- Current AI systems face significant scalability challenges due to reliance on high-quality training data and extensive human workforces for reward design.,   - The proposed "Era of Experience" still depends on extensive human engineering of reward functions, shifting the bottleneck from data curation to reward curation., , -
'@

$ws.Range("L2").Value = @'
- The grounded-agency gap: the inability of contemporary AI systems to autonomously formulate, adapt, and pursue objectives in response to changing circumstances.
'@

$ws.Range("A3").Value = @'
2508.05614v1
'@

$ws.Range("B3").Value = @'
OmniEAR: Benchmarking Agent Reasoning in Embodied Tasks
'@

$ws.Range("C3").Value = @'
Zixuan Wang, Dingming Li, Hongxing Li, Shuo Chen, Yuchen Yan, Wenqi Zhang, Yongliang Shen, Weiming Lu, Jun Xiao, Yueting Zhuang
'@

$ws.Range("D3").Value = @'
http://arxiv.org/abs/2508.05614v1
'@

$ws.Range("E3").Value = @'
  Large language models excel at abstract reasoning but their capacity for
embodied agent reasoning remains largely unexplored. We present OmniEAR, a
comprehensive framework for evaluating how language models reason about
physical interactions, tool usage, and multi-agent coordination in embodied
tasks. Unlike existing benchmarks that provide predefined tool sets or explicit
collaboration directives, OmniEAR requires agents to dynamically acquire
capabilities and autonomously determine coordination strategies based on task
demands. Through text-based environment representation, we model continuous
physical properties and complex spatial relationships across 1,500 scenarios
spanning household and industrial domains. Our systematic evaluation reveals
severe performance degradation when models must reason from constraints: while
achieving 85-96% success with explicit instructions, performance drops to
56-85% for tool reasoning and 63-85% for implicit collaboration, with compound
tasks showing over 50% failure rates. Surprisingly, complete environmental
information degrades coordination performance, indicating models cannot filter
task-relevant constraints. Fine-tuning improves single-agent tasks dramatically
(0.6% to 76.3%) but yields minimal multi-agent gains (1.5% to 5.5%), exposing
fundamental architectural limitations. These findings demonstrate that embodied
reasoning poses fundamentally different challenges than current models can
address, establishing OmniEAR as a rigorous benchmark for evaluating and
advancing embodied AI systems. Our code and data are included in the
supplementary materials and will be open-sourced upon acceptance.

'@

$ws.Range("F3").Value = @'
Project Page: https://zju-real.github.io/OmniEmbodied Code:
  https://github.com/ZJU-REAL/OmniEmbodied
'@

$ws.Range("H3").Value = ""

$ws.Range("I3").Value = @'
- OmniEAR framework,   - Text-based environment representation,   - Fine-tuning of language models,   - Evaluation across 1,500 scenarios spanning household and industrial domains, , -
'@

$ws.Range("J3").Value = @'
- Comprehensive evaluation of language models' reasoning about physical interactions, tool usage, and multi-agent coordination,   - Dynamic acquisition of capabilities and autonomous determination of coordination strategies,   - Systematic evaluation revealing performance insights,   - Establishment of a rigorous benchmark for embodied AI systems,   - Open-sourcing of code and data, , -
'@

$ws.Range("K3").Value = @'
- Severe performance degradation when models must reason from constraints,   - Performance drops significantly for tool reasoning and implicit collaboration tasks,   - High failure rates for compound tasks,   - Degraded coordination performance with complete environmental information,   - Minimal multi-agent gains from fine-tuning, , -
'@

$ws.Range("L3").Value = @'
- The need for improved models that can better handle embodied reasoning tasks,   - The challenge of filtering task-relevant constraints from complete environmental information,   - Fundamental architectural limitations in current models for multi-agent tasks,   - The necessity for further advancements in embodied AI systems to address the unique challenges posed by embodied reasoning
'@

$ws.Range("A4").Value = @'
2508.05519v1
'@

$ws.Range("B4").Value = @'
Leveraging AI to Accelerate Clinical Data Cleaning: A Comparative Study
  of AI-Assisted vs. Traditional Methods
'@

$ws.Range("C4").Value = @'
Matthew Purri, Amit Patel, Erik Deurrell
'@

$ws.Range("D4").Value = @'
http://arxiv.org/abs/2508.05519v1
'@

$ws.Range("E4").Value = @'
  Clinical trial data cleaning represents a critical bottleneck in drug
development, with manual review processes struggling to manage exponentially
increasing data volumes and complexity. This paper presents Octozi, an
artificial intelligence-assisted platform that combines large language models
with domain-specific heuristics to transform clinical data review. In a
controlled experimental study with experienced clinical reviewers (n=10), we
demonstrate that AI assistance increased data cleaning throughput by 6.03-fold
while simultaneously decreasing cleaning errors from 54.67% to 8.48% (a
6.44-fold improvement). Crucially, the system reduced false positive queries by
15.48-fold, minimizing unnecessary site burden. These improvements were
consistent across reviewers regardless of experience level, suggesting broad
applicability. Our findings indicate that AI-assisted approaches can address
fundamental inefficiencies in clinical trial operations, potentially
accelerating drug development timelines and reducing costs while maintaining
regulatory compliance. This work establishes a framework for integrating AI
into safety-critical clinical workflows and demonstrates the transformative
potential of human-AI collaboration in pharmaceutical clinical trials.

'@

$ws.Range("I4").Value = @'
- Artificial Intelligence (AI)-assisted platform: Octozi,   - Large language models,   - Domain-specific heuristics,   - Controlled experimental study with experienced clinical reviewers (n=10), , -
'@

$ws.Range("J4").Value = @'
- Increased data cleaning throughput by 6.03-fold,   - Decreased cleaning errors from 54.67% to 8.48% (a 6.44-fold improvement),   - Reduced false positive queries by 15.48-fold, minimizing unnecessary site burden,   - Consistent improvements across reviewers regardless of experience level,   - Potential to accelerate drug development timelines and reduce costs,   - Maintained regulatory compliance, , -
'@

$ws.Range("K4").Value = @'
- The study's sample size is relatively small (n=10), which may not represent the entire population of clinical reviewers.,   - The abstract does not mention any potential drawbacks or challenges of implementing the AI-assisted platform in real-world settings.,   - The long-term effects and adaptability of the AI system to evolving clinical trial data and protocols are not discussed., , -
'@

$ws.Range("L4").Value = @'
- The abstract does not provide information on how the AI system handles different types of clinical trial data or its adaptability to various trial phases.,   - There is no mention of how the AI system ensures data privacy and security, which is crucial in clinical trials.,   - The abstract does not discuss the potential for bias in the AI system and how it is addressed.,   - Further research is needed to evaluate the system's performance in diverse clinical trial settings and with larger groups of reviewers.
'@

$ws.Range("A5").Value = @'
2508.05469v1
'@

$ws.Range("B5").Value = @'
Let's Measure Information Step-by-Step: LLM-Based Evaluation Beyond
  Vibes
'@

$ws.Range("C5").Value = @'
Zachary Robertson, Sanmi Koyejo
'@

$ws.Range("D5").Value = @'
http://arxiv.org/abs/2508.05469v1
'@

$ws.Range("E5").Value = @'
  We develop mechanisms for evaluating AI systems without ground truth by
exploiting a connection between gaming resistance and output quality. The data
processing inequality ensures post-hoc attempts to game a metric degrades both
information content and task performance. We prove that f-mutual information
measures are the unique gaming resistant mechanisms under natural conditions,
with the overseer acting as an agent. While Shannon mutual information faces
exponential sample complexity, bounded measures like total variation distance
remain tractable. Empirically, across ten domains from translation to peer
review, all information-theoretic mechanisms achieve perfect discrimination (d
> 0.5) between faithful and strategic agents. In contrast, LLM judges exhibit
systematic evaluation inversion, preferring fabricated content over accurate
summaries. Our mechanisms show 10-100x better robustness to adversarial
manipulation than current practices. We also find performance follows an
inverted-U curve with compression ratio, peaking at 10:1 where agent responses
exhibit optimal information diversity (3 effective dimensions), giving a
bias-variance perspective on when our approach is expected to be most
effective.

'@

$ws.Range("F5").Value = @'
13 pages
'@

$ws.Range("A6").Value = @'
2508.05464v1
'@

$ws.Range("B6").Value = @'
Bench-2-CoP: Can We Trust Benchmarking for EU AI Compliance?
'@

$ws.Range("C6").Value = @'
Matteo Prandi, Vincenzo Suriani, Federico Pierucci, Marcello Galisai, Daniele Nardi, Piercosma Bisconti
'@

$ws.Range("D6").Value = @'
http://arxiv.org/abs/2508.05464v1
'@

$ws.Range("E6").Value = @'
  The rapid advancement of General Purpose AI (GPAI) models necessitates robust
evaluation frameworks, especially with emerging regulations like the EU AI Act
and its associated Code of Practice (CoP). Current AI evaluation practices
depend heavily on established benchmarks, but these tools were not designed to
measure the systemic risks that are the focus of the new regulatory landscape.
This research addresses the urgent need to quantify this "benchmark-regulation
gap." We introduce Bench-2-CoP, a novel, systematic framework that uses
validated LLM-as-judge analysis to map the coverage of 194,955 questions from
widely-used benchmarks against the EU AI Act's taxonomy of model capabilities
and propensities. Our findings reveal a profound misalignment: the evaluation
ecosystem is overwhelmingly focused on a narrow set of behavioral propensities,
such as "Tendency to hallucinate" (53.7% of the corpus) and "Discriminatory
bias" (28.9%), while critical functional capabilities are dangerously
neglected. Crucially, capabilities central to loss-of-control scenarios,
including evading human oversight, self-replication, and autonomous AI
development, receive zero coverage in the entire benchmark corpus. This
translates to a near-total evaluation gap for systemic risks like "Loss of
Control" (0.4% coverage) and "Cyber Offence" (0.8% coverage). This study
provides the first comprehensive, quantitative analysis of this gap, offering
critical insights for policymakers to refine the CoP and for developers to
build the next generation of evaluation tools, ultimately fostering safer and
more compliant AI.

'@

$ws.Range("F6").Value = ""

$ws.Range("H6").Value = ""
